$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3988.889
$ws.Range("J76").Value = 5566.6665
$ws.Range("L76").Value = 5566.6665
$ws.Range("N76").Value = -6196.6665
$ws.Range("H79").Value = 3988.889
$ws.Range("J79").Value = 5566.6665
$ws.Range("L79").Value = 5566.6665
$ws.Range("N79").Value = -7750.6665
$ws.Range("H135").Value = 1992.4445
$ws.Range("I135").Value = 1990.2858
$ws.Range("K135").Value = 17912.5722
$ws.Range("M135").Value = -15377.5722
$ws.Range("H138").Value = 3014.8235
$ws.Range("I138").Value = 1248.1
$ws.Range("K138").Value = 3744.3
$ws.Range("M138").Value = 1395.7
$ws.Range("H141").Value = 7781.0435
$ws.Range("I141").Value = 4827.1113
$ws.Range("K141").Value = 14481.3339
$ws.Range("M141").Value = -9301.333899999998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3603.5417
$ws.Range("I2").Value = 3131.2856
$ws.Range("K2").Value = 3131.2856
$ws.Range("M2").Value = -3018.2856
$ws.Range("H88").Value = 4178
$ws.Range("H91").Value = 4178
$ws.Range("H112").Value = 28496.75
$ws.Range("J112").Value = 28496.75
$ws.Range("L112").Value = 28496.75
$ws.Range("N112").Value = -31450.75
$ws.Range("H116").Value = 3603.5417
$ws.Range("I116").Value = 3131.2856
$ws.Range("K116").Value = 3131.2856
$ws.Range("M116").Value = -837.2856000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3603.5417
$ws.Range("I3").Value = 3131.2856
$ws.Range("K3").Value = 3131.2856
$ws.Range("M3").Value = -3017.2856
$ws.Range("H86").Value = 3732.389
$ws.Range("I86").Value = 2304.8823
$ws.Range("J86").Value = 28000
$ws.Range("K86").Value = 2304.8823
$ws.Range("L86").Value = 28000
$ws.Range("M86").Value = -1181.8823
$ws.Range("N86").Value = -30246
$ws.Range("H89").Value = 3732.389
$ws.Range("I89").Value = 2304.8823
$ws.Range("J89").Value = 28000
$ws.Range("K89").Value = 11524.4115
$ws.Range("L89").Value = 140000
$ws.Range("M89").Value = -5908.411500000002
$ws.Range("N89").Value = -151232
$ws.Range("H94").Value = 5347.3076
$ws.Range("I94").Value = 4390.8887
$ws.Range("K94").Value = 4390.8887
$ws.Range("M94").Value = -3939.8887

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3713.6897
$ws.Range("I31").Value = 1675
$ws.Range("J31").Value = 13499.4
$ws.Range("K31").Value = 1675
$ws.Range("L31").Value = 13499.4
$ws.Range("M31").Value = -1380
$ws.Range("N31").Value = -14089.4
$ws.Range("H34").Value = 3713.6897
$ws.Range("I34").Value = 1675
$ws.Range("J34").Value = 13499.4
$ws.Range("K34").Value = 1675
$ws.Range("L34").Value = 13499.4
$ws.Range("M34").Value = -1473
$ws.Range("N34").Value = -13903.4
$ws.Range("H51").Value = 10090
$ws.Range("I51").Value = 10090
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 10090
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9354
$ws.Range("H61").Value = 10090
$ws.Range("I61").Value = 10090
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10090
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -9742
$ws.Range("H96").Value = 13637.25
$ws.Range("J96").Value = 13637.25
$ws.Range("L96").Value = 13637.25
$ws.Range("N96").Value = -19129.25
$ws.Range("H99").Value = 2637.3333
$ws.Range("I99").Value = 2475
$ws.Range("K99").Value = 2475
$ws.Range("M99").Value = -977
$ws.Range("H126").Value = 2637.3333
$ws.Range("I126").Value = 2475
$ws.Range("K126").Value = 7425
$ws.Range("M126").Value = -4955

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2022425.6
$ws.Range("J131").Value = 2224570.5
$ws.Range("L131").Value = 6673711.5
$ws.Range("N131").Value = -6683791.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888
$ws.Range("H70").Value = 14499742
$ws.Range("I70").Value = 30308022
$ws.Range("J70").Value = 8817.083000000001
$ws.Range("K70").Value = 30308022
$ws.Range("L70").Value = 8817.083000000001
$ws.Range("M70").Value = -30307752
$ws.Range("N70").Value = -9357.083000000001
$ws.Range("H73").Value = 14499742
$ws.Range("I73").Value = 30308022
$ws.Range("J73").Value = 8817.083000000001
$ws.Range("K73").Value = 30308022
$ws.Range("L73").Value = 8817.083000000001
$ws.Range("M73").Value = -30307086
$ws.Range("N73").Value = -10689.083
$ws.Range("H97").Value = 1031.25
$ws.Range("I97").Value = 862.4
$ws.Range("J97").Value = 1312.6666
$ws.Range("K97").Value = 862.4
$ws.Range("L97").Value = 1312.6666
$ws.Range("M97").Value = -366.4
$ws.Range("N97").Value = -2304.6666
$ws.Range("H110").Value = 75000
$ws.Range("J110").Value = 75000
$ws.Range("L110").Value = 75000
$ws.Range("N110").Value = -83180
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -46134

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7306.375
$ws.Range("I40").Value = 4655.778
$ws.Range("K40").Value = 4655.778
$ws.Range("M40").Value = -4519.778
$ws.Range("H43").Value = 21499.41
$ws.Range("I43").Value = 19451.62
$ws.Range("K43").Value = 19451.62
$ws.Range("M43").Value = -19258.62
$ws.Range("H55").Value = 1287.125
$ws.Range("I55").Value = 609.8
$ws.Range("J55").Value = 2416
$ws.Range("K55").Value = 609.8
$ws.Range("L55").Value = 2416
$ws.Range("M55").Value = -436.8
$ws.Range("N55").Value = -2762
$ws.Range("H97").Value = 44421.75
$ws.Range("J97").Value = 44421.75
$ws.Range("L97").Value = 44421.75
$ws.Range("N97").Value = -46403.75
$ws.Range("H100").Value = 10244.75
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H132").Value = 3993.8525
$ws.Range("I132").Value = 3524
$ws.Range("K132").Value = 10572
$ws.Range("M132").Value = -8042

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1564.6
$ws.Range("I81").Value = 1734.4
$ws.Range("J81").Value = 1225
$ws.Range("K81").Value = 3468.8
$ws.Range("L81").Value = 2450
$ws.Range("M81").Value = -2407.8
$ws.Range("N81").Value = -4572
$ws.Range("H84").Value = 1564.6
$ws.Range("I84").Value = 1734.4
$ws.Range("J84").Value = 1225
$ws.Range("K84").Value = 17344
$ws.Range("L84").Value = 12250
$ws.Range("M84").Value = -12040
$ws.Range("N84").Value = -22858
$ws.Range("H95").Value = 21385.666
$ws.Range("J95").Value = 21385.666
$ws.Range("L95").Value = 21385.666
$ws.Range("N95").Value = -26877.666
$ws.Range("H99").Value = 37500
$ws.Range("I99").Value = 25000
$ws.Range("K99").Value = 25000
$ws.Range("M99").Value = -22005
$ws.Range("H113").Value = 1792.625
$ws.Range("I113").Value = 1762.2727
$ws.Range("J113").Value = 1859.4
$ws.Range("K113").Value = 5286.8181
$ws.Range("L113").Value = 5578.200000000001
$ws.Range("M113").Value = -3116.8181
$ws.Range("N113").Value = -9918.200000000001
$ws.Range("H132").Value = 3693.5625
$ws.Range("I132").Value = 3946.3704
$ws.Range("J132").Value = 2328.4
$ws.Range("K132").Value = 11839.1112
$ws.Range("L132").Value = 6985.200000000001
$ws.Range("M132").Value = -9309.111199999999
$ws.Range("N132").Value = -12045.2

